# Insert a new row at position 203 (pushing existing rows 203..271 down to 204..272)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with the new record
$ws.Cells.Item(203, 1).Value = 4
$ws.Cells.Item(203, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(203, 3).Value = "Los Lagos"
$ws.Cells.Item(203, 4).Value = 44627
$ws.Cells.Item(203, 5).Value = 10
$ws.Cells.Item(203, 6).Value = 100112045
$ws.Cells.Item(203, 7).Value = "Zapallo"
$ws.Cells.Item(203, 8).Value = "Paine"
$ws.Cells.Item(203, 9).Value = "1a (cosecha)"
$ws.Cells.Item(203, 10).Value = 500
$ws.Cells.Item(203, 11).Value = 500
$ws.Cells.Item(203, 12).Value = 500
$ws.Cells.Item(203, 13).Value = 500
$ws.Cells.Item(203, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(203, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(203, 16).Value = 500
$ws.Cells.Item(203, 17).Value = 1
$ws.Cells.Item(203, 18).Value = "Hortaliza"
